$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: remove the Model (D23) and SN/Lot (E23) values, keep A/B/C as-is
$ws.Range("D23").ClearContents()
$ws.Range("E23").ClearContents()

# New row 24: additional "D2" / "New Item" inventory entry
$ws.Range("A24").Value = "D2"
$ws.Range("B24").Value = "New Item"
$ws.Range("C24").Value = 1

# New row 25: another "D2" / "New Item" inventory entry, with SN/Lot "123"
$ws.Range("A25").Value = "D2"
$ws.Range("B25").Value = "New Item"
$ws.Range("C25").Value = 1
# "123" is a Lot/SN label, not a quantity - force text storage like the
# other SN/Lot cells in this column (e.g. E17/E18) so it round-trips as a
# string instead of a number.
$ws.Range("E25").Value = "'123"
